# Update countries & provincias Spain
# Daily COVID data refresh for the "Pais" sheet:
#  - bump the "Datos actualizados ..." timestamp string (A1)
#  - update Casos totales/Nuevos casos/Casos activos/Recuperados/Casos criticos/Muertes hoy/Muertes
#    (columns B-H) for the countries whose figures moved
#  - three countries changed rank and therefore swapped rows: Namibia now
#    sorts above Guinea Ecuatorial / Republica de Africa Central (rows
#    109-111), and Siria now sorts above Sierra Leona (rows 138-139)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Footer timestamp
$ws.Range("A1").Value = "Datos actualizados a 20 de Agosto de 2020 a las 20:55"

# Estados Unidos (row 4)
$ws.Range("B4").Value = 5721127
$ws.Range("C4").Value = 20196
$ws.Range("D4").Value = 3068577
$ws.Range("E4").Value = 2475749
$ws.Range("G4").Value = 467
$ws.Range("H4").Value = 176801

# India (row 6)
$ws.Range("B6").Value = 2904329
$ws.Range("C6").Value = 68507
$ws.Range("D6").Value = 2157941
$ws.Range("E6").Value = 691413
$ws.Range("G6").Value = 981
$ws.Range("H6").Value = 54975

# España (row 13)
$ws.Range("G13").Value = 16
$ws.Range("H13").Value = 28813

# Argentina (row 21)
$ws.Range("B21").Value = 254520
$ws.Range("C21").Value = 1412
$ws.Range("D21").Value = 234797
$ws.Range("E21").Value = 13665
$ws.Range("G21").Value = 19
$ws.Range("H21").Value = 6058

# Arabia Saudita (row 22)
$ws.Range("B22").Value = 230356
$ws.Range("C22").Value = 656
$ws.Range("E22").Value = 16237
$ws.Range("G22").Value = 5
$ws.Range("H22").Value = 9319

# Canada (row 66)
$ws.Range("B66").Value = 31441
$ws.Range("C66").Value = 426
$ws.Range("D66").Value = 17869
$ws.Range("E66").Value = 13056
$ws.Range("G66").Value = 10
$ws.Range("H66").Value = 516

# Bolivia (row 70)
$ws.Range("B70").Value = 27676
$ws.Range("C70").Value = 129
$ws.Range("E70").Value = 2536
$ws.Range("G70").Value = 1
$ws.Range("H70").Value = 1776

# Republica Dominicana (row 76)
$ws.Range("E76").Value = 7187
$ws.Range("G76").Value = 1
$ws.Range("H76").Value = 120

# Rows 109-111: Namibia jumps above Guinea Ecuatorial / Republica de Africa Central
$ws.Range("A109").Value = "Namibia"
$ws.Range("B109").Value = 4912
$ws.Range("C109").Value = 247
$ws.Range("D109").Value = 2442
$ws.Range("E109").Value = 2429
$ws.Range("G109").Value = 2
$ws.Range("H109").Value = 41

$ws.Range("A110").Value = "Guinea Ecuatorial"
$ws.Range("B110").Value = 4892
$ws.Range("D110").Value = 2713
$ws.Range("E110").Value = 2096
$ws.Range("H110").Value = 83

$ws.Range("A111").Value = "Republica de Africa Central"
$ws.Range("B111").Value = 4679
$ws.Range("D111").Value = 1755
$ws.Range("E111").Value = 2863
$ws.Range("H111").Value = 61

# Hong Kong (row 115)
$ws.Range("B115").Value = 4110
$ws.Range("C115").Value = 52
$ws.Range("D115").Value = 2643
$ws.Range("E115").Value = 1386
$ws.Range("G115").Value = 2
$ws.Range("H115").Value = 81

# Montenegro (row 117)
$ws.Range("B117").Value = 3565
$ws.Range("C117").Value = 83
$ws.Range("D117").Value = 2894
$ws.Range("E117").Value = 583

# Rows 138-139: Siria jumps above Sierra Leona
$ws.Range("A138").Value = "Siria"
$ws.Range("B138").Value = 2008
$ws.Range("C138").Value = 81
$ws.Range("D138").Value = 460
$ws.Range("E138").Value = 1466
$ws.Range("G138").Value = 4
$ws.Range("H138").Value = 82

$ws.Range("A139").Value = "Sierra Leona"
$ws.Range("B139").Value = 1969
$ws.Range("C139").Value = 8
$ws.Range("D139").Value = 1536
$ws.Range("E139").Value = 364
$ws.Range("H139").Value = 69

# Row 161
$ws.Range("B161").Value = 972
$ws.Range("C161").Value = 1
$ws.Range("D161").Value = 869

# Row 190
$ws.Range("B190").Value = 150
$ws.Range("C190").Value = 2
$ws.Range("D190").Value = 115
$ws.Range("E190").Value = 31
